$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column U (21st column), shifting codigo_plaza (old U) to V
$ws.Columns.Item(21).Insert()

# Set the header for the newly inserted column
$ws.Range("U1").Value = "cuspp"
$ws.Range("U1").Style = $ws.Range("T1").Style

# Match the column width that Excel would have computed for the new "cuspp" header
$ws.Columns.Item(21).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Reset selection to default (A1)
$ws.Range("A1").Select()
